# Financials update: add a new "FY2018" (period ending 2018-12-31) column of
# data as the new first data column (D), shifting the existing D:K columns
# (FY2017 .. FY2011) one column to the right (E:L), in each of the three
# statement blocks (Income Statement, Balance Sheet, Cash Flow Statement).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SEB")

# 1) Insert a new blank column across the used data rows, shifting D:K -> E:L
$ws.Range("D5:D102").Insert(-4161)

# 2) New column D should carry the same number formats/styles as the data
#    that was just shifted into column E (date format on header rows, number
#    format on value rows, etc.)
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# 3) Populate the new column D with the FY2018 figures for each block.

# --- Income Statement (rows 7-35), Period Ending 31-Dec-18 ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 6583000
$ws.Range("D9").Value = 6060000
$ws.Range("D10").Value = 523000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 6374000
$ws.Range("D18").Value = 209000
$ws.Range("D20").Value = -181000
$ws.Range("D21").Value = 162000
$ws.Range("D22").Value = 44000
$ws.Range("D23").Value = -16000
$ws.Range("D24").Value = -9000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -7000
$ws.Range("D27").Value = -7000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -10000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 181000
$ws.Range("D33").Value = -17000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -17000

# --- Balance Sheet (rows 38-77), Period Ending 31-Dec-18 ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 194000
$ws.Range("D42").Value = 1336000
$ws.Range("D43").Value = 551000
$ws.Range("D44").Value = 815000
$ws.Range("D45").Value = 131000
$ws.Range("D46").Value = 3027000
$ws.Range("D47").Value = 804000
$ws.Range("D48").Value = 1160000
$ws.Range("D49").Value = 236000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 80000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 5307000
$ws.Range("D57").Value = 238000
$ws.Range("D58").Value = 187000
$ws.Range("D59").Value = 359000
$ws.Range("D60").Value = 784000
$ws.Range("D61").Value = 739000
$ws.Range("D62").Value = 455000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1989000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 3727000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 3318000
$ws.Range("D77").Value = 0

# --- Cash Flow Statement (rows 80-102), Period Ending 31-Dec-18 ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -17000
$ws.Range("D83").Value = 134000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 238000
$ws.Range("D91").Value = -162000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -349000
$ws.Range("D96").Value = -7000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 190000
$ws.Range("D101").Value = -1000
$ws.Range("D102").Value = 78000

Write-Output "SEB financials updated with FY2018 column"
